# Auto-generated Excel COM-interop script applying the Ridill_Profits sheet updates.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 287.37
$ws.Range("I15").Value = 287.37
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 862.11
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -693.11

# Row 132
$ws.Range("H132").Value = 4007704.2
$ws.Range("I132").Value = 837754.1
$ws.Range("J132").Value = 27782330
$ws.Range("K132").Value = 2513262.3
$ws.Range("L132").Value = 83346990
$ws.Range("M132").Value = -2510732.3

# Row 135
$ws.Range("H135").Value = 16667903
$ws.Range("I135").Value = 618.1852
$ws.Range("J135").Value = 166673470
$ws.Range("K135").Value = 5563.6668
$ws.Range("L135").Value = 1500061230
$ws.Range("M135").Value = -3028.6668
$ws.Range("N135").Value = -1500066300

# Row 138
$ws.Range("H138").Value = 3082.47
$ws.Range("I138").Value = 1419.7273
$ws.Range("J138").Value = 3551.4487
$ws.Range("K138").Value = 4259.1819
$ws.Range("L138").Value = 10654.3461
$ws.Range("M138").Value = 880.8181000000004
$ws.Range("N138").Value = -20934.3461

# Row 141
$ws.Range("H141").Value = 3011.2307
$ws.Range("I141").Value = 1008.1667
$ws.Range("J141").Value = 4728.143
$ws.Range("K141").Value = 3024.5001
$ws.Range("L141").Value = 14184.429
$ws.Range("M141").Value = 2155.4999
$ws.Range("N141").Value = -24544.429

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2217749.2
$ws.Range("I32").Value = 2792349.5
$ws.Range("J32").Value = 40316.527
$ws.Range("K32").Value = 2792349.5
$ws.Range("L32").Value = 40316.527
$ws.Range("M32").Value = -2792062.5
$ws.Range("N32").Value = -40890.527

# Row 61
$ws.Range("H61").Value = 2788079
$ws.Range("I61").Value = 1816321.4
$ws.Range("J61").Value = 5349985.5
$ws.Range("K61").Value = 1816321.4
$ws.Range("L61").Value = 5349985.5
$ws.Range("M61").Value = -1816109.4
$ws.Range("N61").Value = -5350409.5

# Row 62
$ws.Range("H62").Value = 18624.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 18624.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 18624.5
$ws.Range("N62").Value = -19872.5

# Row 65
$ws.Range("H65").Value = 18624.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 18624.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 55873.5
$ws.Range("N65").Value = -62113.5

# Row 74
$ws.Range("H74").Value = 18393660
$ws.Range("I74").Value = 1514.5834
$ws.Range("J74").Value = 31376350
$ws.Range("K74").Value = 1514.5834
$ws.Range("L74").Value = 31376350
$ws.Range("M74").Value = -640.5834
$ws.Range("N74").Value = -31378098

# Row 77
$ws.Range("H77").Value = 18393660
$ws.Range("I77").Value = 1514.5834
$ws.Range("J77").Value = 31376350
$ws.Range("K77").Value = 7572.916999999999
$ws.Range("L77").Value = 156881750
$ws.Range("M77").Value = -3204.916999999999
$ws.Range("N77").Value = -156890486

# Row 81
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 10000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -9002
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 10000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -25008
$ws.Range("N84").ClearContents()

# Row 136
$ws.Range("H136").Value = 2788079
$ws.Range("I136").Value = 1816321.4
$ws.Range("J136").Value = 5349985.5
$ws.Range("K136").Value = 5448964.199999999
$ws.Range("L136").Value = 16049956.5
$ws.Range("M136").Value = -5446414.199999999
$ws.Range("N136").Value = -16055056.5

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 9989170
$ws.Range("I134").Value = 11531447
$ws.Range("J134").Value = 74530.71000000001
$ws.Range("K134").Value = 34594341
$ws.Range("L134").Value = 223592.13
$ws.Range("M134").Value = -34591806
$ws.Range("N134").Value = -228662.13

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1814092
$ws.Range("I31").Value = 2526155.5
$ws.Range("J31").Value = 6546
$ws.Range("K31").Value = 2526155.5
$ws.Range("L31").Value = 6546
$ws.Range("M31").Value = -2525860.5
$ws.Range("N31").Value = -7136

# Row 34
$ws.Range("H34").Value = 1814092
$ws.Range("I34").Value = 2526155.5
$ws.Range("J34").Value = 6546
$ws.Range("K34").Value = 2526155.5
$ws.Range("L34").Value = 6546
$ws.Range("M34").Value = -2525953.5
$ws.Range("N34").Value = -6950

# Row 99
$ws.Range("H99").Value = 7838.909
$ws.Range("I99").Value = 11362
$ws.Range("J99").Value = 6517.75
$ws.Range("K99").Value = 11362
$ws.Range("L99").Value = 6517.75
$ws.Range("M99").Value = -9864
$ws.Range("N99").Value = -9513.75

# Row 126
$ws.Range("H126").Value = 7838.909
$ws.Range("I126").Value = 11362
$ws.Range("J126").Value = 6517.75
$ws.Range("K126").Value = 34086
$ws.Range("L126").Value = 19553.25
$ws.Range("M126").Value = -31616
$ws.Range("N126").Value = -24493.25

# Row 132
$ws.Range("H132").Value = 1790.2084
$ws.Range("I132").Value = 1195.125
$ws.Range("J132").Value = 2980.375
$ws.Range("K132").Value = 3585.375
$ws.Range("L132").Value = 8941.125
$ws.Range("M132").Value = -1055.375
$ws.Range("N132").Value = -14001.125

# Row 134
$ws.Range("H134").Value = 801843.9399999999
$ws.Range("I134").Value = 1220.2325
$ws.Range("J134").Value = 5719961
$ws.Range("K134").Value = 3660.6975
$ws.Range("L134").Value = 17159883
$ws.Range("M134").Value = -1125.6975
$ws.Range("N134").Value = -17164953

$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 1200
$ws.Range("I41").Value = 525
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 1575
$ws.Range("L41").Value = 4500
$ws.Range("M41").Value = -1237
$ws.Range("N41").Value = -5176

# Row 113
$ws.Range("H113").Value = 2876.05
$ws.Range("I113").Value = 2556.5
$ws.Range("J113").Value = 3013
$ws.Range("K113").Value = 7669.5
$ws.Range("L113").Value = 9039
$ws.Range("M113").Value = -5499.5
$ws.Range("N113").Value = -13379

# Row 132
$ws.Range("H132").Value = 1732.9375
$ws.Range("I132").Value = 984.3333
$ws.Range("J132").Value = 2025.8695
$ws.Range("K132").Value = 8858.9997
$ws.Range("L132").Value = 18232.8255
$ws.Range("M132").Value = -6328.9997
$ws.Range("N132").Value = -23292.8255

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 8258461.5
$ws.Range("I132").Value = 7504939
$ws.Range("J132").Value = 11366741
$ws.Range("K132").Value = 22514817
$ws.Range("L132").Value = 34100223
$ws.Range("M132").Value = -22512287
$ws.Range("N132").Value = -34105283

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3042
$ws.Range("I40").Value = 3931.3333
$ws.Range("J40").Value = 2375
$ws.Range("K40").Value = 3931.3333
$ws.Range("L40").Value = 2375
$ws.Range("M40").Value = -3795.3333
$ws.Range("N40").Value = -2647

# Row 132
$ws.Range("H132").Value = 3866485.2
$ws.Range("I132").Value = 5106766
$ws.Range("J132").Value = 7833.222
$ws.Range("K132").Value = 15320298
$ws.Range("L132").Value = 23499.666
$ws.Range("M132").Value = -15317768
$ws.Range("N132").Value = -28559.666

# Row 136
$ws.Range("H136").Value = 7814414.5
$ws.Range("I136").Value = 12500865
$ws.Range("J136").Value = 3664.1667
$ws.Range("K136").Value = 37502595
$ws.Range("L136").Value = 10992.5001
$ws.Range("M136").Value = -37500045
$ws.Range("N136").Value = -16092.5001

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 5569.2
$ws.Range("I34").Value = 526
$ws.Range("J34").Value = 6830
$ws.Range("K34").Value = 526
$ws.Range("L34").Value = 6830
$ws.Range("M34").Value = -323
$ws.Range("N34").Value = -7236

